$wb = $excel.ActiveWorkbook
$summary = $wb.Worksheets.Item("Summary")
$lr = $wb.Worksheets.Item("LogisticRegression - Obesity")

# Insert two rows above the current header row so a title row (row1) and a
# blank spacer row (row2) are introduced, pushing the existing table down.
$summary.Rows("1:2").Insert()

# New title cell in A1.
$summary.Range("A1").Value = "Using LR, one hot encoding and ngram(1,2)"

# Re-use the existing bold/red "title" formatting (already used on the
# first sheet) instead of re-building it property by property, which keeps
# the workbook's style table untouched.
$lr.Range("A2").Copy() | Out-Null
$summary.Range("A1").PasteSpecial(-4122) | Out-Null

# Widen column A to fit the new title and disable the old "best fit" autosize.
$summary.Columns("A").ColumnWidth = 38.307291666666664

# Update the remembered selection on the Summary sheet.
$summary.Range("J9").Select() | Out-Null

# Add the new (still empty) third sheet after "Summary".
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $summary)
$newSheet.Name = "Sheet1"
$newSheet.Range("G5").Select() | Out-Null

$summary.Activate() | Out-Null
